$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: dxdate_primary1
$ws.Range("A8").Value = "dxdate_primary1"
$ws.Range("B8").Value = 777
$ws.Range("C8").Value = "valid.changes"
$ws.Range("L8").Value = 7777
$ws.Range("M8").Value = "dxdate_primary1 changed to 7777 from 777 to match data dictionary"

# Row 9: dxdate_primary2
$ws.Range("A9").Value = "dxdate_primary2"
$ws.Range("B9").Value = 777
$ws.Range("C9").Value = "valid.changes"
$ws.Range("L9").Value = 7777
$ws.Range("M9").Value = "dxdate_primary2 changed to 7777 from 777 to match data dictionary"

# Match the shaded/bold-free styling used by the other "changed" rows (A3/A4/A6/A7/M3/M4/M6/M7)
$ws.Range("A3").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("M3").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection down to the next empty row, like Excel does after data entry
$ws.Range("M10").Select() | Out-Null
